# Update "想去人数" (want-to-go count) values for two events
# Sheet "展览" (Exhibition) and "全部类型" (All Types) both contain the
# same underlying data table; update F2 and F7 on each.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 517
    $ws.Range("F7").Value = 725
}
